# Scheduled market-data refresh: updates cached price/profit figures
# on the Leve profit sheets (columns H-N) to the latest API snapshot.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2083.3215
$ws.Range("I17").Value = 1777.7778
$ws.Range("J17").Value = 2228.0527
$ws.Range("K17").Value = 5333.3334
$ws.Range("L17").Value = 6684.158100000001
$ws.Range("M17").Value = -5165.3334
$ws.Range("N17").Value = -7020.158100000001
$ws.Range("H74").Value = 13827
$ws.Range("I74").Value = 3933
$ws.Range("K74").Value = 3933
$ws.Range("M74").Value = -2997
$ws.Range("H77").Value = 13827
$ws.Range("I77").Value = 3933
$ws.Range("K77").Value = 19665
$ws.Range("M77").Value = -14985
$ws.Range("H100").Value = 4710.5
$ws.Range("I100").Value = 2932.875
$ws.Range("K100").Value = 2932.875
$ws.Range("M100").Value = -2391.875
$ws.Range("H116").Value = 2845.1667
$ws.Range("I116").Value = 3265
$ws.Range("K116").Value = 3265
$ws.Range("M116").Value = 177
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = ""
$ws.Range("H135").Value = 2260.5
$ws.Range("I135").Value = 2079.4
$ws.Range("K135").Value = 18714.6
$ws.Range("M135").Value = -16179.6

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4532.6875
$ws.Range("I63").Value = 2633.8333
$ws.Range("K63").Value = 2633.8333
$ws.Range("M63").Value = -1947.8333
$ws.Range("H66").Value = 4532.6875
$ws.Range("I66").Value = 2633.8333
$ws.Range("K66").Value = 13169.1665
$ws.Range("M66").Value = -9737.166499999999
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = ""
$ws.Range("H122").Value = 1112.6
$ws.Range("I122").Value = 1112.6
$ws.Range("K122").Value = 3337.8
$ws.Range("M122").Value = -887.7999999999997

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4638.3335
$ws.Range("I134").Value = 3457.75
$ws.Range("K134").Value = 10373.25
$ws.Range("M134").Value = -7838.25

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4708.5264
$ws.Range("I31").Value = 1766.0435
$ws.Range("J31").Value = 9220.333000000001
$ws.Range("K31").Value = 1766.0435
$ws.Range("L31").Value = 9220.333000000001
$ws.Range("M31").Value = -1471.0435
$ws.Range("N31").Value = -9810.333000000001
$ws.Range("H34").Value = 4708.5264
$ws.Range("I34").Value = 1766.0435
$ws.Range("J34").Value = 9220.333000000001
$ws.Range("K34").Value = 1766.0435
$ws.Range("L34").Value = 9220.333000000001
$ws.Range("M34").Value = -1564.0435
$ws.Range("N34").Value = -9624.333000000001
$ws.Range("H62").Value = 3200
$ws.Range("I62").Value = 3200
$ws.Range("K62").Value = 3200
$ws.Range("M62").Value = -2576
$ws.Range("H65").Value = 3200
$ws.Range("I65").Value = 3200
$ws.Range("K65").Value = 16000
$ws.Range("M65").Value = -12880
$ws.Range("H99").Value = 3014.56
$ws.Range("I99").Value = 2711.5217
$ws.Range("J99").Value = 6499.5
$ws.Range("K99").Value = 2711.5217
$ws.Range("L99").Value = 6499.5
$ws.Range("M99").Value = -1213.5217
$ws.Range("N99").Value = -9495.5
$ws.Range("H126").Value = 3014.56
$ws.Range("I126").Value = 2711.5217
$ws.Range("J126").Value = 6499.5
$ws.Range("K126").Value = 8134.5651
$ws.Range("L126").Value = 19498.5
$ws.Range("M126").Value = -5664.5651
$ws.Range("N126").Value = -24438.5

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 836.375
$ws.Range("I5").Value = 578.6
$ws.Range("J5").Value = 1266
$ws.Range("K5").Value = 1735.8
$ws.Range("L5").Value = 3798
$ws.Range("M5").Value = -1623.8
$ws.Range("N5").Value = -4022
$ws.Range("H131").Value = 300
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = ""
$ws.Range("H135").Value = 836.375
$ws.Range("I135").Value = 578.6
$ws.Range("J135").Value = 1266
$ws.Range("K135").Value = 5207.400000000001
$ws.Range("L135").Value = 11394
$ws.Range("M135").Value = -2672.400000000001
$ws.Range("N135").Value = -16464

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 2250
$ws.Range("I4").Value = 2250
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2250
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""
$ws.Range("H80").Value = 2284.2144
$ws.Range("I80").Value = 2329.8333
$ws.Range("K80").Value = 2329.8333
$ws.Range("M80").Value = -1331.8333
$ws.Range("H83").Value = 2284.2144
$ws.Range("I83").Value = 2329.8333
$ws.Range("K83").Value = 11649.1665
$ws.Range("M83").Value = -6657.166499999999
$ws.Range("H122").Value = 188424.56
$ws.Range("I122").Value = 296864.47
$ws.Range("J122").Value = 4076.7
$ws.Range("K122").Value = 890593.4099999999
$ws.Range("L122").Value = 12230.1
$ws.Range("M122").Value = -888143.4099999999
$ws.Range("N122").Value = -17130.1
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
$ws.Range("H126").Value = 5124.75
$ws.Range("I126").Value = 4999.6665
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 14998.9995
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -12528.9995
$ws.Range("N126").Value = -21440
$ws.Range("H135").Value = 237500
$ws.Range("J135").Value = 237500
$ws.Range("L135").Value = 237500
$ws.Range("N135").Value = -247640

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3493.4285
$ws.Range("I7").Value = 3319.7727
$ws.Range("J7").Value = 4130.1665
$ws.Range("K7").Value = 3319.7727
$ws.Range("L7").Value = 4130.1665
$ws.Range("M7").Value = -3207.7727
$ws.Range("N7").Value = -4354.1665
$ws.Range("H68").Value = 7633.4443
$ws.Range("I68").Value = 2833.3333
$ws.Range("J68").Value = 10033.5
$ws.Range("K68").Value = 2833.3333
$ws.Range("L68").Value = 10033.5
$ws.Range("M68").Value = -2084.3333
$ws.Range("N68").Value = -11531.5
$ws.Range("H71").Value = 7633.4443
$ws.Range("I71").Value = 2833.3333
$ws.Range("J71").Value = 10033.5
$ws.Range("K71").Value = 14166.6665
$ws.Range("L71").Value = 50167.5
$ws.Range("M71").Value = -10422.6665
$ws.Range("N71").Value = -57655.5
$ws.Range("H100").Value = 6428.9375
$ws.Range("I100").Value = 2810.5
$ws.Range("K100").Value = 2810.5
$ws.Range("M100").Value = -2269.5
$ws.Range("H126").Value = 3493.4285
$ws.Range("I126").Value = 3319.7727
$ws.Range("J126").Value = 4130.1665
$ws.Range("K126").Value = 9959.3181
$ws.Range("L126").Value = 12390.4995
$ws.Range("M126").Value = -7489.3181
$ws.Range("N126").Value = -17330.4995
$ws.Range("H132").Value = 3830.682
$ws.Range("I132").Value = 3213.8
$ws.Range("K132").Value = 9641.400000000001
$ws.Range("M132").Value = -7111.400000000001

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1788.5476
$ws.Range("I132").Value = 1686.3846
$ws.Range("K132").Value = 5059.1538
$ws.Range("M132").Value = -2529.1538

